# Append two more years (2021, 2022) of death-total data to the existing
# 2010-2020 series on Sheet1, extending columns L and M, then move the
# active selection to the newly added last cell (M2) to match the
# author's final cursor position.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1: year headers
$ws.Range("L1").Value = 2021
$ws.Range("M1").Value = 2022

# Row 2: yearly death totals
$ws.Range("L2").Value = 81415
$ws.Range("M2").Value = 76633

# Match the saved selection/active cell from the source edit
$null = $ws.Range("M2").Select()
